# "haciendo los ajustes finales al juego" - add the two missing players
# (Ruben, Jacob) to the game-stats table and leave the selection where
# the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook no longer carries (unused) protection settings
$wb.Unprotect()

# Row 8 - Ruben: 1 partida jugada, 1 perdida, 0% de aciertos
$ws.Range("A8").Value = "Ruben"
$ws.Range("B8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0

# Row 9 - Jacob: 1 partida jugada, 1 ganada, 100% de aciertos
$ws.Range("A9").Value = "Jacob"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("E9").Value = 100

# Restore the cursor/selection to where it was when the file was saved
$null = $ws.Range("H19").Select()
